$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case municipality / state names (Excel PROPER-style) ---
$ws.Range('B8').Value = 'Pabellón De Arteaga'
$ws.Range('B9').Value = 'Rincón De Romos'
$ws.Range('B10').Value = 'San Francisco De Los Romo'
$ws.Range('B11').Value = 'San José De Gracia'
$ws.Range('B29').Value = 'Amatenango De La Frontera'
$ws.Range('B32').Value = 'Benemérito De Las Américas'
$ws.Range('B37').Value = 'Comitán De Domínguez'
$ws.Range('B56').Value = 'Salto De Agua'
$ws.Range('B57').Value = 'San Cristóbal De Las Casas'
$ws.Range('B92').Value = 'Coyame Del Sotol'
$ws.Range('B103').Value = 'Guadalupe Y Calvo'
$ws.Range('B106').Value = 'Hidalgo Del Parral'
$ws.Range('B128').Value = 'San Francisco De Borja'
$ws.Range('B129').Value = 'San Francisco Del Oro'
$ws.Range('B137').Value = 'Valle De Zaragoza'
$ws.Range('B151').Value = 'San Juan De Sabinas'
$ws.Range('B162').Value = 'Villa De Álvarez'
$ws.Range('A164').Value = 'Ciudad De México'
$ws.Range('B168').Value = 'Cuajimalpa De Morelos'
$ws.Range('B183').Value = 'Coneto De Comonfort'
$ws.Range('B197').Value = 'Nombre De Dios'
$ws.Range('B200').Value = 'Pánuco De Coronado'
$ws.Range('B207').Value = 'San Juan De Guadalupe'
$ws.Range('B208').Value = 'San Juan Del Río'
$ws.Range('B209').Value = 'San Luis Del Cordero'
$ws.Range('B210').Value = 'San Pedro Del Gallo'
$ws.Range('A220').Value = 'Estado De México'
$ws.Range('B220').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B223').Value = 'Almoloya De Alquisiras'
$ws.Range('B224').Value = 'Almoloya De Juárez'
$ws.Range('B227').Value = 'Atizapán De Zaragoza'
$ws.Range('B238').Value = 'Ecatepec De Morelos'
$ws.Range('B241').Value = 'Ixtapan De La Sal'
$ws.Range('B249').Value = 'Naucalpan De Juárez'
$ws.Range('B255').Value = 'San Felipe Del Progreso'
$ws.Range('B256').Value = 'San Martín De Las Pirámides'
$ws.Range('B258').Value = 'San Simón De Guerrero'
$ws.Range('B267').Value = 'Tenango Del Valle'
$ws.Range('B274').Value = 'Tlalnepantla De Baz'
$ws.Range('B277').Value = 'Valle De Bravo'
$ws.Range('B278').Value = 'Villa De Allende'
$ws.Range('B279').Value = 'Villa Del Carbón'
$ws.Range('B295').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B299').Value = 'Jaral Del Progreso'
$ws.Range('B306').Value = 'Purísima Del Rincón'
$ws.Range('B311').Value = 'San Francisco Del Rincón'
$ws.Range('B313').Value = 'San Luis De La Paz'
$ws.Range('B315').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B316').Value = 'Silao De La Victoria'
$ws.Range('B318').Value = 'Valle De Santiago'
$ws.Range('B324').Value = 'Acapulco De Juárez'
$ws.Range('B326').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B331').Value = 'Atoyac De Álvarez'
$ws.Range('B332').Value = 'Ayutla De Los Libres'
$ws.Range('B335').Value = 'Chilapa De Álvarez'
$ws.Range('B336').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B337').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B341').Value = 'Coyuca De Benítez'
$ws.Range('B342').Value = 'Coyuca De Catalán'
$ws.Range('B346').Value = 'Cuetzala Del Progreso'
$ws.Range('B347').Value = 'Cutzamala De Pinzón'
$ws.Range('B352').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B353').Value = 'Iguala De La Independencia'
$ws.Range('B355').Value = 'Zihuatanejo De Azueta'
$ws.Range('B357').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B371').Value = 'Taxco De Alarcón'
$ws.Range('B373').Value = 'Técpan De Galeana'
$ws.Range('B376').Value = 'Tixtla De Guerrero'
$ws.Range('B379').Value = 'Tlapa De Comonfort'
$ws.Range('B393').Value = 'Atotonilco El Grande'
$ws.Range('B398').Value = 'Cuautepec De Hinojosa'
$ws.Range('B405').Value = 'Jacala De Ledezma'
$ws.Range('B410').Value = 'Mixquiahuala De Juárez'
$ws.Range('B411').Value = 'Omitlán De Juárez'
$ws.Range('B412').Value = 'Pachuca De Soto'
$ws.Range('B413').Value = 'Progreso De Obregón'
$ws.Range('B421').Value = 'Tepehuacán De Guerrero'
$ws.Range('B422').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B424').Value = 'Tezontepec De Aldama'
$ws.Range('B430').Value = 'Tula De Allende'
$ws.Range('B431').Value = 'Tulancingo De Bravo'
$ws.Range('B433').Value = 'Zacualtipán De Ángeles'
$ws.Range('B437').Value = 'Ahualulco De Mercado'
$ws.Range('B440').Value = 'Atemajac De Brizuela'
$ws.Range('B443').Value = 'Atotonilco El Alto'
$ws.Range('B445').Value = 'Autlán De Navarro'
$ws.Range('B454').Value = 'Cuautitlán De García Barragán'
$ws.Range('B461').Value = 'Encarnación De Díaz'
$ws.Range('B466').Value = 'Huejuquilla El Alto'
$ws.Range('B467').Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range('B468').Value = 'Ixtlahuacán Del Río'
$ws.Range('B471').Value = 'Jilotlán De Los Dolores'
$ws.Range('B476').Value = 'Lagos De Moreno'
$ws.Range('B483').Value = 'Ojuelos De Jalisco'
$ws.Range('B486').Value = 'San Diego De Alejandría'
$ws.Range('B487').Value = 'San Juan De Los Lagos'
$ws.Range('B488').Value = 'San Juanito De Escobedo'
$ws.Range('B490').Value = 'San Martín De Bolaños'
$ws.Range('B492').Value = 'San Miguel El Alto'
$ws.Range('B493').Value = 'Santa María De Los Ángeles'
$ws.Range('B495').Value = 'Talpa De Allende'
$ws.Range('B496').Value = 'Tamazula De Gordiano'
$ws.Range('B502').Value = 'Teocuitatlán De Corona'
$ws.Range('B503').Value = 'Tepatitlán De Morelos'
$ws.Range('B505').Value = 'Tizapán El Alto'
$ws.Range('B506').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B513').Value = 'Unión De San Antonio'
$ws.Range('B514').Value = 'Unión De Tula'
$ws.Range('B515').Value = 'Valle De Juárez'
$ws.Range('B520').Value = 'Yahualica De González Gallo'
$ws.Range('B521').Value = 'Zacoalco De Torres'
$ws.Range('B523').Value = 'Zapotlán El Grande'
$ws.Range('B540').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B542').Value = 'Cojumatlán De Régules'
$ws.Range('B590').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B619').Value = 'Puente De Ixtla'
$ws.Range('B623').Value = 'Tetela Del Volcán'
$ws.Range('B625').Value = 'Tlaltizapán De Zapata'
$ws.Range('B631').Value = 'Zacualpan De Amilpas'
$ws.Range('B635').Value = 'Amatlán De Cañas'
$ws.Range('B636').Value = 'Bahía De Banderas'
$ws.Range('B640').Value = 'Ixtlán Del Río'
$ws.Range('B654').Value = 'Ciénega De Flores'
$ws.Range('B659').Value = 'Mier Y Noriega'
$ws.Range('B662').Value = 'San Nicolás De Los Garza'
$ws.Range('B665').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B667').Value = 'Guevea De Humboldt'
$ws.Range('B668').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B669').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B670').Value = 'Ixtlán De Juárez'
$ws.Range('B671').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B673').Value = 'Mariscala De Juárez'
$ws.Range('B675').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B678').Value = 'Nejapa De Madero'
$ws.Range('B679').Value = 'Oaxaca De Juárez'
$ws.Range('B680').Value = 'Ocotlán De Morelos'
$ws.Range('B681').Value = 'Pinotepa De Don Luis'
$ws.Range('B682').Value = 'Putla Villa De Guerrero'
$ws.Range('B696').Value = 'San Juan Del Estado'
$ws.Range('B709').Value = 'San Miguel Del Puerto'
$ws.Range('B719').Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range('B725').Value = 'Santa Cruz Tacache De Mina'
$ws.Range('B733').Value = 'Santiago Del Río'
$ws.Range('B744').Value = 'Santo Domingo De Morelos'
$ws.Range('B747').Value = 'Sitio De Xitlapehua'
$ws.Range('B748').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B750').Value = 'Teotitlán De Flores Magón'
$ws.Range('B751').Value = 'Tezoatlán De Segura Y Luna'
$ws.Range('B752').Value = 'Tlacolula De Matamoros'
$ws.Range('B754').Value = 'Villa De Etla'
$ws.Range('B755').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B756').Value = 'Villa De Zaachila'
$ws.Range('B757').Value = 'Villa Sola De Vega'
$ws.Range('B758').Value = 'Zapotitlán Del Río'
$ws.Range('B759').Value = 'Zimatlán De Álvarez'
$ws.Range('B774').Value = 'Chalchicomula De Sesma'
$ws.Range('B787').Value = 'Huehuetlán El Chico'
$ws.Range('B790').Value = 'Ixcamilpa De Guerrero'
$ws.Range('B793').Value = 'Izúcar De Matamoros'
$ws.Range('B800').Value = 'Los Reyes De Juárez'
$ws.Range('B804').Value = 'Palmar De Bravo'
$ws.Range('B813').Value = 'San Nicolás De Los Ranchos'
$ws.Range('B815').Value = 'San Salvador El Seco'
$ws.Range('B817').Value = 'Tecali De Herrera'
$ws.Range('B822').Value = 'Tepanco De López'
$ws.Range('B823').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B825').Value = 'Tepexi De Rodríguez'
$ws.Range('B826').Value = 'Tetela De Ocampo'
$ws.Range('B831').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B846').Value = 'Amealco De Bonfil'
$ws.Range('B847').Value = 'Cadereyta De Montes'
$ws.Range('B851').Value = 'Jalpan De Serra'
$ws.Range('B852').Value = 'Landa De Matamoros'
$ws.Range('B855').Value = 'Pinal De Amoles'
$ws.Range('B857').Value = 'San Juan Del Río'
$ws.Range('B868').Value = 'Ciudad Del Maíz'
$ws.Range('B876').Value = 'Mexquitic De Carmona'
$ws.Range('B882').Value = 'Santa María Del Río'
$ws.Range('B886').Value = 'Villa De Guadalupe'
$ws.Range('B887').Value = 'Villa De Ramos'
$ws.Range('B888').Value = 'Villa De Reyes'
$ws.Range('B924').Value = 'Nacozari De García'
$ws.Range('B937').Value = 'Jalpa De Méndez'
$ws.Range('B958').Value = 'Soto La Marina'
$ws.Range('B968').Value = 'Nanacamilpa De Mariano Arista'
$ws.Range('B970').Value = 'Papalotla De Xicohténcatl'
$ws.Range('B971').Value = 'Tepetitla De Lardizábal'
$ws.Range('B981').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B984').Value = 'Amatlán De Los Reyes'
$ws.Range('B992').Value = 'Castillo De Teayo'
$ws.Range('B1000').Value = 'Cosamaloapan De Carpio'
$ws.Range('B1010').Value = 'Hueyapan De Ocampo'
$ws.Range('B1011').Value = 'Ignacio De La Llave'
$ws.Range('B1014').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B1018').Value = 'Juchique De Ferrer'
$ws.Range('B1021').Value = 'Las Vigas De Ramírez'
$ws.Range('B1022').Value = 'Lerdo De Tejada'
$ws.Range('B1024').Value = 'Martínez De La Torre'
$ws.Range('B1025').Value = 'Medellín De Bravo'
$ws.Range('B1035').Value = 'Paso Del Macho'
$ws.Range('B1039').Value = 'Poza Rica De Hidalgo'
$ws.Range('B1044').Value = 'Sayula De Alemán'
$ws.Range('B1046').Value = 'Soledad De Doblado'
$ws.Range('B1048').Value = 'Tatahuicapan De Juárez'
$ws.Range('B1061').Value = 'Vega De Alatorre'
$ws.Range('B1067').Value = 'Zozocolco De Hidalgo'
$ws.Range('B1085').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B1095').Value = 'Jiménez Del Teul'
$ws.Range('B1104').Value = 'Nochistlán De Mejía'
$ws.Range('B1105').Value = 'Noria De Ángeles'
$ws.Range('B1115').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B1118').Value = 'Villa De Cos'

# --- Recomputed percentage values (1-ulp float changes) ---
$ws.Range('D97').Value = 0.0009588749200937568
$ws.Range('D103').Value = 0.0009588749200937568
$ws.Range('D157').Value = 0.0009588749200937568
$ws.Range('D176').Value = 0.0009588749200937568
$ws.Range('D177').Value = 0.0009588749200937568
$ws.Range('D179').Value = 0.0009588749200937568
$ws.Range('D190').Value = 0.0009588749200937568
$ws.Range('D191').Value = 0.0009588749200937568
$ws.Range('D195').Value = 0.0009588749200937568
$ws.Range('D198').Value = 0.0009588749200937568
$ws.Range('D228').Value = 0.0009588749200937568
$ws.Range('D244').Value = 0.0009588749200937568
$ws.Range('D306').Value = 0.0009588749200937568
$ws.Range('D309').Value = 0.0009588749200937568
$ws.Range('D339').Value = 0.0009588749200937568
$ws.Range('D387').Value = 0.0009588749200937568
$ws.Range('D447').Value = 0.0009588749200937568
$ws.Range('D496').Value = 0.0009588749200937568
$ws.Range('D523').Value = 0.0009588749200937568
$ws.Range('D551').Value = 0.0009588749200937568
$ws.Range('D679').Value = 0.0009588749200937568
$ws.Range('D686').Value = 0.0009588749200937568
$ws.Range('D698').Value = 0.0009588749200937568
$ws.Range('D770').Value = 0.0009588749200937568
$ws.Range('D867').Value = 0.0009588749200937568
$ws.Range('D896').Value = 0.0009588749200937568
$ws.Range('D910').Value = 0.0009588749200937568
$ws.Range('D951').Value = 0.0009588749200937568
$ws.Range('D1016').Value = 0.0009588749200937568
$ws.Range('D1070').Value = 0.0009588749200937568
$ws.Range('D1118').Value = 0.0009588749200937568
$ws.Range("D1124").Value = 0.09588749200937569

# --- Remove trailing footnote rows (1126-1131) ---
$ws.Range("A1126:D1131").EntireRow.Delete()
